# deploy: path to startup
# Update the deployment version name and the startup script path,
# then move the saved view/selection state to match the new location
# the author was working at.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("application deploy commands")

# Bump the deployment "version name" used throughout the formulas in
# column D (they all reference $B$4, so they recalc automatically).
$ws.Range("B4").Value = "20190928-01"

# The startup script now lives under the active release folder rather
# than directly under /usr/local/casualapp.
$ws.Range("D52").Value = "sh /usr/local/casualapp/active/startup.sh"

# Force a full recalculation so the dependent formula cells (D13, D15,
# D17, D18, D21, D24, D25, D30, D33, D36, D38, D39, ...) pick up the
# new version name.
$excel.CalculateFullRebuild()

# Restore the author's view/selection position after editing.
$ws.Activate()
$ws.Range("D53").Select()
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1
